$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1,1).Value = "Language"
$ws.Cells.Item(1,2).Value = "German"
$ws.Cells.Item(1,3).Value = "Spanish"
$ws.Cells.Item(1,4).Value = "Slovak"
$ws.Cells.Item(1,5).Value = "Norwegian"
$ws.Cells.Item(1,6).Value = "Greek"
$ws.Cells.Item(1,7).Value = "Chinese"
$ws.Cells.Item(1,8).Value = "Vietnamese"
$ws.Cells.Item(1,9).Value = "Indonesian"
$ws.Cells.Item(1,10).Value = "Finnish"
$ws.Cells.Item(1,11).Value = "Basque"
$ws.Cells.Item(1,12).Value = "Korean"
$ws.Cells.Item(1,13).Value = "Japanese"
$ws.Cells.Item(1,14).Value = "Turkish"
$ws.Cells.Item(1,15).Value = "Arabic"
$ws.Cells.Item(1,16).Value = "Hebrew"
$ws.Cells.Item(2,1).Value = "German"
$ws.Cells.Item(2,2).Value = 0.9787187739463602
$ws.Cells.Item(2,3).Value = 0.81036015325670496
$ws.Cells.Item(2,4).Value = 0.73371647509578541
$ws.Cells.Item(2,5).Value = 0.85489348659003828
$ws.Cells.Item(2,6).Value = 0.80654099616858232
$ws.Cells.Item(2,7).Value = 0.64673103448275859
$ws.Cells.Item(2,8).Value = 0.55739157088122604
$ws.Cells.Item(2,9).Value = 0.70211800766283528
$ws.Cells.Item(2,10).Value = 0.68326130268199237
$ws.Cells.Item(2,11).Value = 0.62637854406130267
$ws.Cells.Item(2,12).Value = 0.58102988505747122
$ws.Cells.Item(2,13).Value = 0.53822835249042145
$ws.Cells.Item(2,14).Value = 0.67488429118773952
$ws.Cells.Item(2,15).Value = 0.64718467432950189
$ws.Cells.Item(2,16).Value = 0.73836628352490419
$ws.Cells.Item(3,1).Value = "Spanish"
$ws.Cells.Item(3,2).Value = 0.82568511960000757
$ws.Cells.Item(3,3).Value = 0.98691312664532871
$ws.Cells.Item(3,4).Value = 0.86333592166815021
$ws.Cells.Item(3,5).Value = 0.88649836177345132
$ws.Cells.Item(3,6).Value = 0.78441696179996589
$ws.Cells.Item(3,7).Value = 0.5815041381792011
$ws.Cells.Item(3,8).Value = 0.58214806537754971
$ws.Cells.Item(3,9).Value = 0.83271150167610464
$ws.Cells.Item(3,10).Value = 0.62604874907672203
$ws.Cells.Item(3,11).Value = 0.58966686237003085
$ws.Cells.Item(3,12).Value = 0.52040681047707427
$ws.Cells.Item(3,13).Value = 0.49495274710706239
$ws.Cells.Item(3,14).Value = 0.62426847976364086
$ws.Cells.Item(3,15).Value = 0.65163538569345281
$ws.Cells.Item(3,16).Value = 0.75441752997102329
$ws.Cells.Item(4,1).Value = "Slovak"
$ws.Cells.Item(4,2).Value = 0.8269112680380718
$ws.Cells.Item(4,3).Value = 0.8875498925391464
$ws.Cells.Item(4,4).Value = 0.96369358305188824
$ws.Cells.Item(4,5).Value = 0.85600245624808102
$ws.Cells.Item(4,6).Value = 0.81716303346638008
$ws.Cells.Item(4,7).Value = 0.68360454405894999
$ws.Cells.Item(4,8).Value = 0.6137549892539147
$ws.Cells.Item(4,9).Value = 0.76197420939514893
$ws.Cells.Item(4,10).Value = 0.78661344795824373
$ws.Cells.Item(4,11).Value = 0.74623887012588275
$ws.Cells.Item(4,12).Value = 0.61805342339576297
$ws.Cells.Item(4,13).Value = 0.63286766963463315
$ws.Cells.Item(4,14).Value = 0.74117285845870429
$ws.Cells.Item(4,15).Value = 0.7724132637396377
$ws.Cells.Item(4,16).Value = 0.7645072152287381
$ws.Cells.Item(5,1).Value = "Norwegian"
$ws.Cells.Item(5,2).Value = 0.86748314756724287
$ws.Cells.Item(5,3).Value = 0.85733831675899352
$ws.Cells.Item(5,4).Value = 0.78512313955816593
$ws.Cells.Item(5,5).Value = 0.97707401721951548
$ws.Cells.Item(5,6).Value = 0.80654742040979777
$ws.Cells.Item(5,7).Value = 0.68143896415938066
$ws.Cells.Item(5,8).Value = 0.57485149836481342
$ws.Cells.Item(5,9).Value = 0.78972835880664749
$ws.Cells.Item(5,10).Value = 0.74901555095775207
$ws.Cells.Item(5,11).Value = 0.6452980044049923
$ws.Cells.Item(5,12).Value = 0.58272709070279649
$ws.Cells.Item(5,13).Value = 0.52332643662817857
$ws.Cells.Item(5,14).Value = 0.66281786024160716
$ws.Cells.Item(5,15).Value = 0.69198424881532405
$ws.Cells.Item(5,16).Value = 0.78849362610959084
$ws.Cells.Item(6,1).Value = "Greek"
$ws.Cells.Item(6,2).Value = 0.82219373741073065
$ws.Cells.Item(6,3).Value = 0.82054568760300306
$ws.Cells.Item(6,4).Value = 0.77815418421534521
$ws.Cells.Item(6,5).Value = 0.77952755905511806
$ws.Cells.Item(6,6).Value = 0.97491301959348109
$ws.Cells.Item(6,7).Value = 0.58221937374107302
$ws.Cells.Item(6,8).Value = 0.49661234206189342
$ws.Cells.Item(6,9).Value = 0.69181468595495332
$ws.Cells.Item(6,10).Value = 0.67057315509979853
$ws.Cells.Item(6,11).Value = 0.60355246291887932
$ws.Cells.Item(6,12).Value = 0.54513825306720376
$ws.Cells.Item(6,13).Value = 0.49560520051272661
$ws.Cells.Item(6,14).Value = 0.6399011170115364
$ws.Cells.Item(6,15).Value = 0.66535433070866146
$ws.Cells.Item(6,16).Value = 0.72074711591283647
$ws.Cells.Item(7,1).Value = "Chinese"
$ws.Cells.Item(7,2).Value = 0.53488178488178484
$ws.Cells.Item(7,3).Value = 0.4408924408924409
$ws.Cells.Item(7,4).Value = 0.36771561771561773
$ws.Cells.Item(7,5).Value = 0.62054612054612057
$ws.Cells.Item(7,6).Value = 0.61038961038961037
$ws.Cells.Item(7,7).Value = 0.96386946386946382
$ws.Cells.Item(7,8).Value = 0.52822177822177818
$ws.Cells.Item(7,9).Value = 0.57101232101232102
$ws.Cells.Item(7,10).Value = 0.59948384948384947
$ws.Cells.Item(7,11).Value = 0.5815850815850816
$ws.Cells.Item(7,12).Value = 0.64244089244089242
$ws.Cells.Item(7,13).Value = 0.65068265068265063
$ws.Cells.Item(7,14).Value = 0.52314352314352319
$ws.Cells.Item(7,15).Value = 0.50924075924075929
$ws.Cells.Item(7,16).Value = 0.66441891441891443
$ws.Cells.Item(8,1).Value = "Vietnamese"
$ws.Cells.Item(8,2).Value = 0.45211208699289002
$ws.Cells.Item(8,3).Value = 0.60025094102885823
$ws.Cells.Item(8,4).Value = 0.61639481388540363
$ws.Cells.Item(8,5).Value = 0.61672940192388126
$ws.Cells.Item(8,6).Value = 0.58996235884567128
$ws.Cells.Item(8,7).Value = 0.58912588874947724
$ws.Cells.Item(8,8).Value = 0.89401923881221246
$ws.Cells.Item(8,9).Value = 0.6311166875784191
$ws.Cells.Item(8,10).Value = 0.5692179004600586
$ws.Cells.Item(8,11).Value = 0.59155165202843996
$ws.Cells.Item(8,12).Value = 0.52061898787118366
$ws.Cells.Item(8,13).Value = 0.52354663320786277
$ws.Cells.Item(8,14).Value = 0.55918025930572979
$ws.Cells.Item(8,15).Value = 0.58954412379757426
$ws.Cells.Item(8,16).Value = 0.63429527394395646
$ws.Cells.Item(9,1).Value = "Thai"
$ws.Cells.Item(9,2).Value = 0.57252934324881288
$ws.Cells.Item(9,3).Value = 0.5664366992205
$ws.Cells.Item(9,4).Value = 0.54506764626825555
$ws.Cells.Item(9,5).Value = 0.59224083863453092
$ws.Cells.Item(9,6).Value = 0.61123555236985938
$ws.Cells.Item(9,7).Value = 0.67117641788370219
$ws.Cells.Item(9,8).Value = 0.57436609622793655
$ws.Cells.Item(9,9).Value = 0.62248006451034854
$ws.Cells.Item(9,10).Value = 0.6097123913627811
$ws.Cells.Item(9,11).Value = 0.56778066481498068
$ws.Cells.Item(9,12).Value = 0.51814353552549053
$ws.Cells.Item(9,13).Value = 0.54703879580682735
$ws.Cells.Item(9,14).Value = 0.47401666517337149
$ws.Cells.Item(9,15).Value = 0.59385359734790788
$ws.Cells.Item(9,16).Value = 0.70038527013708451
$ws.Cells.Item(10,1).Value = "Indonesian"
$ws.Cells.Item(10,2).Value = 0.68641765704584046
$ws.Cells.Item(10,3).Value = 0.83786078098471983
$ws.Cells.Item(10,4).Value = 0.77478777589134129
$ws.Cells.Item(10,5).Value = 0.83684210526315794
$ws.Cells.Item(10,6).Value = 0.69736842105263153
$ws.Cells.Item(10,7).Value = 0.68743633276740235
$ws.Cells.Item(10,8).Value = 0.59533106960950766
$ws.Cells.Item(10,9).Value = 0.9286078098471986
$ws.Cells.Item(10,10).Value = 0.74966044142614596
$ws.Cells.Item(10,11).Value = 0.69719864176570456
$ws.Cells.Item(10,12).Value = 0.61663837011884548
$ws.Cells.Item(10,13).Value = 0.58225806451612905
$ws.Cells.Item(10,14).Value = 0.67292020373514427
$ws.Cells.Item(10,15).Value = 0.61315789473684212
$ws.Cells.Item(10,16).Value = 0.77164685908319186
$ws.Cells.Item(11,1).Value = "Finnish"
$ws.Cells.Item(11,2).Value = 0.77695730379627004
$ws.Cells.Item(11,3).Value = 0.85141531761810096
$ws.Cells.Item(11,4).Value = 0.76521821452238947
$ws.Cells.Item(11,5).Value = 0.83072990627662591
$ws.Cells.Item(11,6).Value = 0.83058790116444192
$ws.Cells.Item(11,7).Value = 0.74069866515194549
$ws.Cells.Item(11,8).Value = 0.59736817192085578
$ws.Cells.Item(11,9).Value = 0.76261478746568212
$ws.Cells.Item(11,10).Value = 0.96757549938464449
$ws.Cells.Item(11,11).Value = 0.80199753857805545
$ws.Cells.Item(11,12).Value = 0.67745905519265359
$ws.Cells.Item(11,13).Value = 0.72190665530625764
$ws.Cells.Item(11,14).Value = 0.7824008330966582
$ws.Cells.Item(11,15).Value = 0.67130550033134528
$ws.Cells.Item(11,16).Value = 0.84379437659755752
$ws.Cells.Item(12,1).Value = "Basque"
$ws.Cells.Item(12,2).Value = 0.69885944038729797
$ws.Cells.Item(12,3).Value = 0.69643882825962089
$ws.Cells.Item(12,4).Value = 0.68872569131041272
$ws.Cells.Item(12,5).Value = 0.65918601788791331
$ws.Cells.Item(12,6).Value = 0.672027570361861
$ws.Cells.Item(12,7).Value = 0.6639862148190695
$ws.Cells.Item(12,8).Value = 0.5811110199392796
$ws.Cells.Item(12,9).Value = 0.69570033642405837
$ws.Cells.Item(12,10).Value = 0.72400919012062037
$ws.Cells.Item(12,11).Value = 0.94563879543776153
$ws.Cells.Item(12,12).Value = 0.6292770985476327
$ws.Cells.Item(12,13).Value = 0.71477804217608931
$ws.Cells.Item(12,14).Value = 0.70698285057848531
$ws.Cells.Item(12,15).Value = 0.57618774103552961
$ws.Cells.Item(12,16).Value = 0.73045048001969315
$ws.Cells.Item(13,1).Value = "Korean"
$ws.Cells.Item(13,2).Value = 0.56190509765211871
$ws.Cells.Item(13,3).Value = 0.58961432701121064
$ws.Cells.Item(13,4).Value = 0.52012973277867869
$ws.Cells.Item(13,5).Value = 0.56715786504970744
$ws.Cells.Item(13,6).Value = 0.56980187548473527
$ws.Cells.Item(13,7).Value = 0.57198759077769157
$ws.Cells.Item(13,8).Value = 0.50433617711344569
$ws.Cells.Item(13,9).Value = 0.56870901783825711
$ws.Cells.Item(13,10).Value = 0.58538391031516601
$ws.Cells.Item(13,11).Value = 0.5846435873933582
$ws.Cells.Item(13,12).Value = 0.95409997884791653
$ws.Cells.Item(13,13).Value = 0.61693576817316509
$ws.Cells.Item(13,14).Value = 0.57322146231403792
$ws.Cells.Item(13,15).Value = 0.49104561799337237
$ws.Cells.Item(13,16).Value = 0.59116547979976031
$ws.Cells.Item(14,1).Value = "Japanese"
$ws.Cells.Item(14,2).Value = 0.39855923576853808
$ws.Cells.Item(14,3).Value = 0.30561428235846838
$ws.Cells.Item(14,4).Value = 0.26622817320491737
$ws.Cells.Item(14,5).Value = 0.41288857567927328
$ws.Cells.Item(14,6).Value = 0.46049643724062328
$ws.Cells.Item(14,7).Value = 0.53895544593219014
$ws.Cells.Item(14,8).Value = 0.39942056221125988
$ws.Cells.Item(14,9).Value = 0.35345705113146969
$ws.Cells.Item(14,10).Value = 0.46253229974160209
$ws.Cells.Item(14,11).Value = 0.53699788583509511
$ws.Cells.Item(14,12).Value = 0.58656330749354002
$ws.Cells.Item(14,13).Value = 0.97917156056690935
$ws.Cells.Item(14,14).Value = 0.53597995458460579
$ws.Cells.Item(14,15).Value = 0.35557121603633229
$ws.Cells.Item(14,16).Value = 0.47787957090282668
$ws.Cells.Item(15,1).Value = "Turkish"
$ws.Cells.Item(15,2).Value = 0.59469475326701771
$ws.Cells.Item(15,3).Value = 0.70196996294129121
$ws.Cells.Item(15,4).Value = 0.65837721864638188
$ws.Cells.Item(15,5).Value = 0.67797932514140824
$ws.Cells.Item(15,6).Value = 0.70031207333723422
$ws.Cells.Item(15,7).Value = 0.69582601911449193
$ws.Cells.Item(15,8).Value = 0.64784474351472598
$ws.Cells.Item(15,9).Value = 0.68938950653403552
$ws.Cells.Item(15,10).Value = 0.72683830700214547
$ws.Cells.Item(15,11).Value = 0.74292958845328649
$ws.Cells.Item(15,12).Value = 0.61858786814901501
$ws.Cells.Item(15,13).Value = 0.67749171055197976
$ws.Cells.Item(15,14).Value = 0.93251414082309347
$ws.Cells.Item(15,15).Value = 0.58074897600936215
$ws.Cells.Item(15,16).Value = 0.70392042129900523
$ws.Cells.Item(16,1).Value = "Arabic"
$ws.Cells.Item(16,2).Value = 0.59759265958496399
$ws.Cells.Item(16,3).Value = 0.70135166244614722
$ws.Cells.Item(16,4).Value = 0.68471075739139009
$ws.Cells.Item(16,5).Value = 0.66941822606636625
$ws.Cells.Item(16,6).Value = 0.68148781530568614
$ws.Cells.Item(16,7).Value = 0.49988489492551058
$ws.Cells.Item(16,8).Value = 0.54648600651165846
$ws.Cells.Item(16,9).Value = 0.62465879567204918
$ws.Cells.Item(16,10).Value = 0.54549939158746341
$ws.Cells.Item(16,11).Value = 0.53306804354260529
$ws.Cells.Item(16,12).Value = 0.46390633735652981
$ws.Cells.Item(16,13).Value = 0.43355148485546091
$ws.Cells.Item(16,14).Value = 0.54004012234025056
$ws.Cells.Item(16,15).Value = 0.97128950570592298
$ws.Cells.Item(16,16).Value = 0.80705100799158092
$ws.Cells.Item(17,1).Value = "Hebrew"
$ws.Cells.Item(17,2).Value = 0.55096782717843695
$ws.Cells.Item(17,3).Value = 0.64537226663143288
$ws.Cells.Item(17,4).Value = 0.54449362489264719
$ws.Cells.Item(17,5).Value = 0.5448239413357997
$ws.Cells.Item(17,6).Value = 0.65911343066657857
$ws.Cells.Item(17,7).Value = 0.51674704366783375
$ws.Cells.Item(17,8).Value = 0.4899253484838475
$ws.Cells.Item(17,9).Value = 0.62357138138336521
$ws.Cells.Item(17,10).Value = 0.51575609433837621
$ws.Cells.Item(17,11).Value = 0.47743938693268151
$ws.Cells.Item(17,12).Value = 0.46118781792957653
$ws.Cells.Item(17,13).Value = 0.42775979388253949
$ws.Cells.Item(17,14).Value = 0.48893439915438991
$ws.Cells.Item(17,15).Value = 0.67853603752394798
$ws.Cells.Item(17,16).Value = 0.97245160864107816

$ws.Range("A18:XFD19").Select() | Out-Null